# Marlborough_Sauvignon_blanc_PhenologyObs.xlsx update
# - calculate fraction of bud burst at the budding phase, updating the met and pheno data set
#
# Concretely:
#   1. Rename "Sheet1" -> "Data"
#   2. Rename the .met site filenames (append site codes BRA/OYB/SEA/RPC/VLA)
#   3. Clear the "Text" number-format override on A2:A63 (keep it on the A1 header)
#   4. Resize/declare columns A-C ready for the new fraction-of-budburst columns
#   5. Reset the sheet view: scroll back to the top, select C3 instead of the old G71

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet tab.
$ws.Name = "Data"

# 2. Update the shared .met site-name strings used throughout column A.
$ws.Cells.Replace("PhenoTestClimateSite1001.met", "PhenoTestClimateSite1001_BRA.met") | Out-Null
$ws.Cells.Replace("PhenoTestClimateSite1002.met", "PhenoTestClimateSite1002_OYB.met") | Out-Null
$ws.Cells.Replace("PhenoTestClimateSite1003.met", "PhenoTestClimateSite1003_SEA.met") | Out-Null
$ws.Cells.Replace("PhenoTestClimateSite1004.met", "PhenoTestClimateSite1004_RPC.met") | Out-Null
$ws.Cells.Replace("PhenoTestClimateSite1005.met", "PhenoTestClimateSite1005_VLA.met") | Out-Null

# 3. The data rows (A2:A63) no longer need the explicit "Text" style - only the
#    header (A1) keeps it.
$ws.Range("A2:A63").Style = "Normal"

# 4. Widen column A and declare columns B and C so the new fraction-of-budburst
#    data has room (bestFit-style widths).
$ws.Columns.Item(1).ColumnWidth = 30.666666666666668
$ws.Columns.Item(2).ColumnWidth = 9
$ws.Columns.Item(3).ColumnWidth = 36.666666666666664

# 5. Scroll back to the top-left and select C3 (was scrolled to A47 / G71 selected).
$ws.Range("A1").Select() | Out-Null
$ws.Range("C3").Select() | Out-Null
